$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row numbers (column G, "Recorded By") whose value needs to flip from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System".
$rows = @(
    3,6,10,12,13,14,15,18,19,20,21,22,24,26,29,32,36,38,39,40,41,
    44,45,46,47,48,50,52,55,58,62,64,65,66,67,70,71,72,73,74,76,78,
    83,84,85,86,90,92,99,101,109,110,111,112,116,118,125,127,135,
    136,137,138,142,144,151,153
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
